$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 13.7

$ws.Range("B3").Value = -15.7
$ws.Range("C3").Value = 6.3

$ws.Range("C4").Value = 1.2

$ws.Range("C6").Value = -25.9

$ws.Range("C7").Value = -19.1

$ws.Range("C9").Value = -23.3

$ws.Range("C13").Value = 13.5

$ws.Range("C15").Value = 39.2

$ws.Range("C16").Value = 24.1

$ws.Range("C23").Value = 17.3

$ws.Range("C24").Value = 18.1
